$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.62%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.41%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.995"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.59%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07820"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.75%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.213"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.75%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.030"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.07%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.028"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.51%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9153"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.54%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09728"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.24%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1886"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.41%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08698"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.85%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03564"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.26%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09966"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.62%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001482"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.08%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005641"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.39%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.37%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "7.24%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3462"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.34%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.08%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.763"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.57%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.94%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04618"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.08%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.14%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004787"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.15%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.91%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "38.89%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01772"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.73%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04750"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.34%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008053"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.16%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.34%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007658"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "11.26%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002160"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.89%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009865"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.06%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006013"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.41%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.27%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.907"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "189.37%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.55%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.27%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.27%"
